{"js": "// Update the date line and every two-digit-by-two-digit multiplication\n// problem in the table with its new value. Each \"old\" string is unique in\n// the document, so a body.search() + insertText(\"Replace\") per pair is\n// safe and cannot cascade into an already-updated cell.\n\nconst replacements = [\n  [\"2024-05-28 Tuesday\", \"2024-05-29 Wednesday\"],\n  [\"84\u00d776=\", \"49\u00d788=\"],\n  [\"11\u00d783=\", \"43\u00d792=\"],\n  [\"42\u00d723=\", \"16\u00d773=\"],\n  [\"73\u00d781=\", \"39\u00d750=\"],\n  [\"45\u00d784=\", \"40\u00d722=\"],\n  [\"73\u00d787=\", \"37\u00d798=\"],\n  [\"18\u00d714=\", \"87\u00d795=\"],\n  [\"46\u00d771=\", \"32\u00d718=\"],\n  [\"95\u00d742=\", \"23\u00d780=\"],\n  [\"69\u00d796=\", \"51\u00d724=\"],\n  [\"11\u00d762=\", \"25\u00d772=\"],\n  [\"94\u00d739=\", \"82\u00d761=\"],\n  [\"23\u00d760=\", \"95\u00d754=\"],\n  [\"77\u00d717=\", \"96\u00d754=\"],\n  [\"19\u00d765=\", \"99\u00d719=\"],\n  [\"91\u00d797=\", \"59\u00d735=\"],\n  [\"90\u00d775=\", \"68\u00d717=\"],\n  [\"56\u00d782=\", \"84\u00d722=\"],\n  [\"75\u00d797=\", \"28\u00d794=\"],\n  [\"18\u00d783=\", \"20\u00d766=\"],\n  [\"59\u00d714=\", \"59\u00d769=\"],\n  [\"52\u00d754=\", \"58\u00d775=\"],\n  [\"16\u00d796=\", \"28\u00d742=\"],\n  [\"85\u00d753=\", \"31\u00d717=\"],\n  [\"49\u00d722=\", \"64\u00d744=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and replace each two-digit-by-two-digit multiplication\n# problem with its new value. Every \"old\" string below is unique in the\n# document, so a plain Find/Replace (ReplaceAll) for each pair is safe and\n# will not cascade into an already-replaced cell.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-05-28 Tuesday\", \"2024-05-29 Wednesday\"),\n    @(\"84\u00d776=\", \"49\u00d788=\"),\n    @(\"11\u00d783=\", \"43\u00d792=\"),\n    @(\"42\u00d723=\", \"16\u00d773=\"),\n    @(\"73\u00d781=\", \"39\u00d750=\"),\n    @(\"45\u00d784=\", \"40\u00d722=\"),\n    @(\"73\u00d787=\", \"37\u00d798=\"),\n    @(\"18\u00d714=\", \"87\u00d795=\"),\n    @(\"46\u00d771=\", \"32\u00d718=\"),\n    @(\"95\u00d742=\", \"23\u00d780=\"),\n    @(\"69\u00d796=\", \"51\u00d724=\"),\n    @(\"11\u00d762=\", \"25\u00d772=\"),\n    @(\"94\u00d739=\", \"82\u00d761=\"),\n    @(\"23\u00d760=\", \"95\u00d754=\"),\n    @(\"77\u00d717=\", \"96\u00d754=\"),\n    @(\"19\u00d765=\", \"99\u00d719=\"),\n    @(\"91\u00d797=\", \"59\u00d735=\"),\n    @(\"90\u00d775=\", \"68\u00d717=\"),\n    @(\"56\u00d782=\", \"84\u00d722=\"),\n    @(\"75\u00d797=\", \"28\u00d794=\"),\n    @(\"18\u00d783=\", \"20\u00d766=\"),\n    @(\"59\u00d714=\", \"59\u00d769=\"),\n    @(\"52\u00d754=\", \"58\u00d775=\"),\n    @(\"16\u00d796=\", \"28\u00d742=\"),\n    @(\"85\u00d753=\", \"31\u00d717=\"),\n    @(\"49\u00d722=\", \"64\u00d744=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
